# Apply the changes described by the diff: add "carrier" (D) values to the
# practice rows (2-5) and generic word rows (6-9), add "pair_kind" (J)
# values of unique_video/unique_audio to rows 6-9, and populate rows 14-21
# (the former blank filler rows) with kind (C) and carrier (D) values for
# the new unique_video / unique_audio entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Practice rows (2-5): add column D (carrier) ---
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# --- Generic word rows (6-9): add column J (pair_kind) ---
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# --- Rows 14-17: unique_video kind/carrier rows ---
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"

# --- Rows 18-21: unique_audio kind/carrier rows ---
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
